# update evidence format (A1-A6)
#
# For the "Info" sheet and sheets "A1".."A6", the worksheet used to contain:
#   row 1: column headers
#   row 2: a human-readable description / placeholder ("team name", "addr1", ...)
#   row 3 (and, for "A2", also row 4): example/sample evidence data
#
# The new format drops the description row entirely: the example data moves
# up to become row 2 (and row 3 for "A2", which has two example rows), while
# keeping the formatting that used to belong to row 2 (the description row).

$wb = $excel.ActiveWorkbook

function Update-EvidenceSheet {
    param([string]$SheetName, [string]$LastCol, [int]$DataRows, [string]$SelectCell)

    $ws = $wb.Worksheets.Item($SheetName)

    # Old layout: row 1 = headers, row 2 = description (styled), rows
    # 3..(2+DataRows) = example data (unstyled).
    #
    # Overwrite row 2's VALUES with the first example row's values, while
    # keeping row 2's own cell formatting (it stays the styled row).  Then
    # delete old row 3 outright: everything further down (any additional
    # example rows) shifts up by one in a single step, carrying its own
    # original formatting with it - no extra copy needed for those rows.
    $destRange = "A2:" + $LastCol + "2"
    $srcRange = "A3:" + $LastCol + "3"
    $ws.Range($destRange).Value2 = $ws.Range($srcRange).Value2
    $ws.Rows("3:3").Delete()

    # Park the selection where it ended up after the edit.
    $ws.Range($SelectCell).Select()
}

# Info sheet: TeamName..Community (A:H), one example row. Saved cursor
# lands one row below the data (H3), past the now-2-row used range.
Update-EvidenceSheet "Info" "H" 1 "H3"

# A1: TxHash / ClassID (A:B), one example row.
Update-EvidenceSheet "A1" "B" 1 "B2"

# A2: TxHash / ClassID / NFTID (A:C), two example rows (NFT mint evidence).
Update-EvidenceSheet "A2" "C" 2 "C3"

# A3: TxHash / ClassID / NFTID / ChainID (A:D), one example row.
Update-EvidenceSheet "A3" "D" 1 "D2"

# A4: same shape as A3.
Update-EvidenceSheet "A4" "D" 1 "D2"

# A5: same shape as A3.
Update-EvidenceSheet "A5" "D" 1 "D2"

# A6: same shape as A3.
Update-EvidenceSheet "A6" "D" 1 "D2"

# Restore "A6" as the active sheet/tab, matching the workbook's saved
# activeTab and that sheet's own tabSelected view state.
$wb.Worksheets.Item("A6").Select()
